$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Fix the title heading: "Participant Information" -> "Risk Management Plan"
# ---------------------------------------------------------------------
$headingRange = $d.Content
$headingRange.Find.Execute("Participant Information") | Out-Null
$headingRange.Expand(4) | Out-Null

$headingXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body>' + `
    '<w:p w14:paraId="595B6BAE" w14:textId="77777777" w:rsidR="001F0505" w:rsidRPr="00C53DF7" w:rsidRDefault="001F0505" w:rsidP="001F0505" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
    '<w:pPr><w:rPr><w:b/></w:rPr></w:pPr>' + `
    '<w:r w:rsidRPr="00C53DF7"><w:rPr><w:b/></w:rPr><w:t>Risk Management Plan</w:t></w:r>' + `
    '</w:p>' + `
    '</w:body></w:document>' + `
    '</pkg:xmlData></pkg:part></pkg:package>'

$headingRange.InsertXML($headingXml) | Out-Null

# ---------------------------------------------------------------------
# 2) Convert the ERGO/<Faculty fldSimple>/17661 field from a simple field
#    (<w:fldSimple>) into the equivalent complex field
#    (begin / instrText / separate / result / end run sequence).
# ---------------------------------------------------------------------
$fld = $d.Fields.Item(1)
$fldCode = $fld.Code.Text
$fldResult = $fld.Result.Text

$tbl = $d.Tables.Item(1)
$cell = $tbl.Cell(1, 1)
$fieldPara = $cell.Range.Paragraphs.Item(1)
$fieldParaRange = $fieldPara.Range

$fieldXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body>' + `
    '<w:p w14:paraId="233EDA4A" w14:textId="77777777" w:rsidR="001F0505" w:rsidRPr="00C026B7" w:rsidRDefault="001F0505" w:rsidP="00733384" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
    '<w:pPr><w:spacing w:before="120" w:after="120"/></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">Ethics reference number:  </w:t></w:r>' + `
    '<w:r w:rsidRPr="00B0488E"><w:rPr><w:b/></w:rPr><w:t>ERGO/</w:t></w:r>' + `
    '<w:r><w:fldChar w:fldCharType="begin"/></w:r>' + `
    '<w:r><w:instrText xml:space="preserve">' + $fldCode + '</w:instrText></w:r>' + `
    '<w:r><w:fldChar w:fldCharType="separate"/></w:r>' + `
    '<w:r><w:t>' + $fldResult + '</w:t></w:r>' + `
    '<w:r><w:fldChar w:fldCharType="end"/></w:r>' + `
    '<w:r w:rsidRPr="00B0488E"><w:rPr><w:b/></w:rPr><w:t>/</w:t></w:r>' + `
    '<w:r w:rsidRPr="00E9688A"><w:rPr><w:b/></w:rPr><w:t>17661</w:t></w:r>' + `
    '</w:p>' + `
    '</w:body></w:document>' + `
    '</pkg:xmlData></pkg:part></pkg:package>'

$fieldParaRange.InsertXML($fieldXml) | Out-Null

# ---------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark: it now sits right after "In this study, pa"
#    instead of right after "which provi" (the surrounding run text is
#    re-split accordingly, while the visible paragraph text is unchanged).
# ---------------------------------------------------------------------
$introRange = $d.Content
$introRange.Find.Execute("In this study, participants will be required") | Out-Null
$introRange.Expand(4) | Out-Null

$introXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body>' + `
    '<w:p w14:paraId="7D60716C" w14:textId="1855D934" w:rsidR="002B213F" w:rsidRDefault="0061695A" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
    '<w:r><w:t>In this study, pa</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
    '<w:bookmarkEnd w:id="0"/>' + `
    '<w:r><w:t xml:space="preserve">rticipants will be required to wear a device on various parts of their body which will measure their movement while they perform different exercises. These activities will take place inside a university </w:t></w:r>' + `
    '<w:r w:rsidR="00507072"><w:t xml:space="preserve">building </w:t></w:r>' + `
    '<w:r><w:t>which provides sufficient space, most likely in the level 3 Zepler labs. Due to the nature of this study, there are some risks whose details and management strategies are listed in the table below.</w:t></w:r>' + `
    '</w:p>' + `
    '</w:body></w:document>' + `
    '</pkg:xmlData></pkg:part></pkg:package>'

$introRange.InsertXML($introXml) | Out-Null

Write-Output "edit complete"
